# Uncheck "Snap to grid" (the "If defined, snap to grid when document grid
# is defined" paragraph option) in the "Normal" style.
#
# In the Word object model this checkbox is exposed as
# ParagraphFormat.DisableLineHeightGrid; turning it on writes the OOXML
# paragraph property <w:snapToGrid w:val="0"/> into the style's <w:pPr>.
$d = $word.ActiveDocument

$normal = $d.Styles("Normal")
$normal.ParagraphFormat.DisableLineHeightGrid = $true
